$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.362.20'
$ws.Range("D3").Value = '3.751.40'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.32'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.90'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").Value = '3.750.04'
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("E10").Value = '  -3.39%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("E13").Value = '  -7.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.04'
$ws.Range("D15").Value = '4.379.49'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").Value = '3.756.12'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '68.326.84'
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("E18").Value = '  -4.13%  '
$ws.Range("E19").Value = '  -2.67%  '
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '467.63'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("E23").Value = '  -3.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.50'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000144'
$ws.Range("E25").Value = '  -2.68%  '
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.97'
$ws.Range("E27").Value = '  -1.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '3.897.39'
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  -4.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.29'
$ws.Range("E32").Value = '  -4.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.86'
$ws.Range("E33").Value = '  -2.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.18'
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.21'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D37").Value = '3.706.38'
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("E38").Value = '  -2.78%  '
$ws.Range("E39").Value = '  -11.06%  '
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.60'
$ws.Range("E46").Value = '  -0.83%  '
$ws.Range("E47").Value = '  -1.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.83'
$ws.Range("E48").Value = '  +9.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.83'
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '146.71'
$ws.Range("E50").Value = '  +4.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '391.38'
$ws.Range("E51").Value = '  -1.35%  '
